$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-16"

# Update the shared header string (column I1 label) to the new date
$ws.Range("I1").Value = "2022 (through 10-16)"

# Update the changed data values (October, November rows and Total row)
$ws.Range("I10").Value = 145
$ws.Range("I11").Value = 53
$ws.Range("I14").Value = 1330
